$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top, pushing the existing table down from
# rows 1-12 to rows 3-14.
$ws.Rows("1:2").Insert() | Out-Null

# Add the new title row in A1 (bold, matching the sheet's existing bold style).
$ws.Range("A1").Value = "357 Portus Delphini Ct - Nobili"
$ws.Range("A1").Font.Bold = $true

# Page setup: landscape orientation, 88% scale, fit-to-page enabled.
$ws.PageSetup.Zoom = 88
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1
$ws.PageSetup.Orientation = 2

# Move the active selection to A15, just below the last data row.
$ws.Range("A15").Select() | Out-Null
